# Regenerate save_data to use K (strikeouts) instead of the old Strike# derived
# value. The new K values are sourced from the authoritative per-game data
# (not derivable from the other columns already on the sheet), so we write
# the literal values back into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 4
    25 = 3
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 2
    32 = 3
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 3
    39 = 0
    40 = 1
    41 = 2
    42 = 1
    43 = 0
    44 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
